$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value for every data row
# (rows 2-497). The commit updates that date from 45186 (2023-09-17)
# to 45188 (2023-09-19) for every row.
$ws.Range("C2:C497").Value = 45188
